# RAD Production Release Nov 08 2024 - refresh the recorded run timestamps
# in column B (rows 2-30) with the latest Katalon test-run log values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$timestamps = @(
    "Tue Nov 12 17:29:13 EST 2024",
    "Tue Nov 12 17:29:25 EST 2024",
    "Tue Nov 12 17:29:35 EST 2024",
    "Tue Nov 12 17:29:45 EST 2024",
    "Tue Nov 12 17:29:55 EST 2024",
    "Tue Nov 12 17:30:05 EST 2024",
    "Tue Nov 12 17:30:15 EST 2024",
    "Tue Nov 12 17:30:25 EST 2024",
    "Tue Nov 12 17:30:36 EST 2024",
    "Tue Nov 12 17:30:49 EST 2024",
    "Tue Nov 12 17:30:59 EST 2024",
    "Tue Nov 12 17:31:09 EST 2024",
    "Tue Nov 12 17:31:19 EST 2024",
    "Tue Nov 12 17:31:30 EST 2024",
    "Tue Nov 12 17:31:41 EST 2024",
    "Tue Nov 12 17:31:51 EST 2024",
    "Tue Nov 12 17:32:01 EST 2024",
    "Tue Nov 12 17:32:11 EST 2024",
    "Tue Nov 12 17:32:22 EST 2024",
    "Tue Nov 12 17:32:32 EST 2024",
    "Tue Nov 12 17:32:42 EST 2024",
    "Tue Nov 12 17:32:52 EST 2024",
    "Tue Nov 12 17:33:02 EST 2024",
    "Tue Nov 12 17:33:12 EST 2024",
    "Tue Nov 12 17:33:22 EST 2024",
    "Tue Nov 12 17:33:32 EST 2024",
    "Tue Nov 12 17:33:44 EST 2024",
    "Tue Nov 12 17:33:55 EST 2024",
    "Tue Nov 12 17:34:05 EST 2024"
)

# Rows 2..30 in column B get the new timestamp strings, in order.
for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Range("B$row").Value = $timestamps[$i]
}
